$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-03-26"

# Update the header label in I1 to match the new date
$ws.Range("I1").Value = "2022 (through 03-26)"

# Update the data values that changed for 2022-04-03
$ws.Range("I4").Value = 104
$ws.Range("I14").Value = 404
